$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns keep their original "text" semantics.
# Many of the values look like plain numbers (e.g. "4.81", "91.06") and
# Excel would otherwise silently convert them to numeric cells, losing the
# original formatting (and, for rows like "35.648.79" which use multiple
# dots, that particular value is safe -- but to be consistent and avoid any
# accidental numeric coercion we force column D and E to Text format before
# writing any values).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "35.711.58"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "1.986.37"
$ws.Range("E3").Value = "  -3.95%  "
$ws.Range("D5").Value = "242.44"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "0.639"
$ws.Range("E6").Value = "  -4.31%  "
$ws.Range("D7").Value = "57.14"
$ws.Range("E7").Value = "  +8.13%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "59.66"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").Value = "0.361"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "0.0731"
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("E12").Value = "  -4.68%  "
$ws.Range("D13").Value = "0.926"
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("D14").Value = "14.13"
$ws.Range("E14").Value = "  -2.70%  "
$ws.Range("D15").Value = "2.275.27"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("D17").Value = "1.985.48"
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("E18").Value = "  +5.41%  "
$ws.Range("D19").Value = "35.539.35"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "70.55"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").Value = "233.11"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").Value = "5.07"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").Value = "2.34"
$ws.Range("E26").Value = "  +10.10%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.16"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "163.57"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  -5.34%  "
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").Value = "4.81"
$ws.Range("E32").Value = "  -5.68%  "
$ws.Range("D33").Value = "0.0589"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").Value = "0.0898"
$ws.Range("E34").Value = "  +9.66%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "2.37"
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "4.27"
$ws.Range("E36").Value = "  -6.58%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("E42").Value = "  -2.63%  "
$ws.Range("E43").Value = "  -4.16%  "
$ws.Range("D44").Value = "0.0892"
$ws.Range("E44").Value = "  -4.69%  "
$ws.Range("D45").Value = "91.06"
$ws.Range("D46").Value = "1.383.41"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "7.42"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").Value = "15.45"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "2.88"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("D51").Value = "45.79"
$ws.Range("E51").Value = "  +1.96%  "
